$d = $word.ActiveDocument

# Locate the paragraph that contains the long "tests" blurb with the two
# Laravel course hyperlinks (rId6-rId9). We find it by its distinctive
# leading text rather than a hard-coded paragraph index so the script
# stays robust if earlier content shifts.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Les tests sont un sujet complexe*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq 0) {
    Write-Output "target paragraph not found"
} else {
    $target = $d.Paragraphs($targetIndex)

    # Remove the whole paragraph (text, the two hyperlink runs, and its
    # trailing paragraph mark) in one shot.
    $target.Range.Delete()

    # The paragraph immediately after (already blank, matching the
    # surrounding style) has now shifted into this slot. Clone it twice by
    # inserting fresh blank paragraphs right before it - this reproduces
    # the "split into two empty paragraphs" shape from the target edit.
    $after = $d.Paragraphs($targetIndex)
    $after.Range.InsertParagraphBefore()
    $after2 = $d.Paragraphs($targetIndex + 1)
    $after2.Range.InsertParagraphBefore()
}
